$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for changed rows.
# NumberFormat is forced to text ("@") before writing D-column values so
# that numeric-looking strings (e.g. "572.16", "0.999") are preserved as
# literal text instead of being parsed into floating point numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.011.35"
$ws.Range("E2").Value = "  +1.68%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.164.28"
$ws.Range("E3").Value = "  +3.45%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.16"
$ws.Range("E5").Value = "  +2.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.53"
$ws.Range("E6").Value = "  +5.56%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.162.67"
$ws.Range("E8").Value = "  +3.41%  "

$ws.Range("E9").Value = "  +2.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  +4.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.19"
$ws.Range("E11").Value = "  +2.24%  "

$ws.Range("E12").Value = "  +4.84%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000274"
$ws.Range("E13").Value = "  +18.70%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.05"
$ws.Range("E14").Value = "  +7.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.679.96"
$ws.Range("E15").Value = "  +3.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.096.98"
$ws.Range("E16").Value = "  +1.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.164.88"
$ws.Range("E17").Value = "  +3.50%  "

$ws.Range("E18").Value = "  +6.33%  "

$ws.Range("E19").Value = "  +1.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "509.93"
$ws.Range("E20").Value = "  +6.78%  "

$ws.Range("E21").Value = "  +6.67%  "

$ws.Range("E24").Value = "  +2.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.70"
$ws.Range("E25").Value = "  +3.20%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.16"
$ws.Range("E27").Value = "  +14.93%  "

$ws.Range("E28").Value = "  +3.68%  "

$ws.Range("E29").Value = "  +8.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.05"
$ws.Range("E30").Value = "  +6.65%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.80"
$ws.Range("E31").Value = "  +15.40%  "

$ws.Range("E32").Value = "  +7.31%  "

$ws.Range("E33").Value = "  +0.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.32"
$ws.Range("E34").Value = "  +11.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.66"
$ws.Range("E35").Value = "  +6.94%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.61"
$ws.Range("E36").Value = "  +1.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "478.53"
$ws.Range("E37").Value = "  +7.93%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0885"
$ws.Range("E38").Value = "  +9.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.09"
$ws.Range("E39").Value = "  +9.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0420"
$ws.Range("E40").Value = "  +3.23%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.122.71"
$ws.Range("E41").Value = "  +4.33%  "

$ws.Range("E42").Value = "  +4.48%  "

$ws.Range("E43").Value = "  +4.97%  "

$ws.Range("E44").Value = "  +17.59%  "

$ws.Range("E45").Value = "  +10.79%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.32"
$ws.Range("E46").Value = "  +5.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0594"
$ws.Range("E47").Value = "  +14.46%  "

$ws.Range("E49").Value = "  +1.75%  "

$ws.Range("E50").Value = "  +11.30%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "122.48"
$ws.Range("E51").Value = "  +2.78%  "

# Rows 22 and 23 swap coin identities (Polygon <-> InternetComputer(DFINITY))
# along with their corresponding link/price/volume values.
$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.733"
$ws.Range("E22").Value = "  +7.58%  "

$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.64"
$ws.Range("E23").Value = "  +8.48%  "

